$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "... (to be used as a final test), 500 to 1000 names becoming the ..."
#    -> "... (to be used as a final test), the following 500 names becoming the ..."
#    "500" must stay bold+italic (matching the old "500 to 1000" run's formatting).
# ---------------------------------------------------------------------------

# 1a. turn the plain ", " right after "final test)" into ", the following "
$r1 = $d.Content
$found1 = $r1.Find.Execute("(to be used as a final test), ", $true, $false, $false, $false, $false, $true, 1, $false, "(to be used as a final test), the following ", 2)
Write-Host "step1a found: $found1"

# 1b. shrink the bold+italic "500 to 1000" run down to "500 " (keeps its formatting)
$r2 = $d.Content
$found2 = $r2.Find.Execute("500 to 1000", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "step1b found: $found2"
$r2.Text = "500 "
$r2.Collapse(0)

# 1c. the plain run after the bold one used to start with a space (" names becoming");
#     that space is now redundant since "500 " already ends with one. Keep this
#     replacement inside the plain run only, so the following bold+underline
#     "dev-test set" run keeps its formatting untouched.
$r3 = $d.Content
$found3 = $r3.Find.Execute(" names becoming", $true, $false, $false, $false, $false, $true, 1, $false, "names becoming", 2)
Write-Host "step1c found: $found3"

# 1d. Word re-anchors the (invisible) _GoBack bookmark to the last edit location;
#     move it to sit right after the shrunk "500 " run, matching that behaviour.
$d.Bookmarks("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $r2)

# ---------------------------------------------------------------------------
# 2) "... last letter of the names was printed out. Using ..." (unchanged text,
#    merges "was " + "print" + "ed" + " out" into one run) - no visible change.
# ---------------------------------------------------------------------------
$r4 = $d.Content
$found4 = $r4.Find.Execute(" was printed out", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "step2 found: $found4"
$r4.Text = " was printed outZZZTMP"
$r4b = $d.Content
$found4b = $r4b.Find.Execute(" was printed outZZZTMP", $true, $false, $false, $false, $false, $true, 1, $false, " was printed out", 2)
Write-Host "step2b found: $found4b"

# ---------------------------------------------------------------------------
# 3) '... last letter being "e" would allow ...' -> '... last letter being "e", would allow ...'
# ---------------------------------------------------------------------------
$r5 = $d.Content
$found5 = $r5.Find.Execute("last letter being ""e"" would allow", $true, $false, $false, $false, $false, $true, 1, $false, "last letter being ""e"", would allow", 2)
Write-Host "step3 found: $found5"

# ---------------------------------------------------------------------------
# 4) '... is female but "Eddie" and yet ...' -> '... is female but "Eddie" is male and yet ...'
# ---------------------------------------------------------------------------
$r6 = $d.Content
$found6 = $r6.Find.Execute(" ""Eddie"" and yet", $true, $false, $false, $false, $false, $true, 1, $false, " ""Eddie"" is male and yet", 2)
Write-Host "step4 found: $found6"

# ---------------------------------------------------------------------------
# 5) "Due combining" -> "Because combining"
# ---------------------------------------------------------------------------
$r7 = $d.Content
$found7 = $r7.Find.Execute("Due combining", $true, $false, $false, $false, $false, $true, 1, $false, "Because combining", 2)
Write-Host "step5 found: $found7"

# ---------------------------------------------------------------------------
# 6) "first the letters and the first four letters" -> "first the 3 letters and the first 4 letters"
# ---------------------------------------------------------------------------
$r8 = $d.Content
$found8 = $r8.Find.Execute("first the letters and the first four letters", $true, $false, $false, $false, $false, $true, 1, $false, "first the 3 letters and the first 4 letters", 2)
Write-Host "step6 found: $found8"

# ---------------------------------------------------------------------------
# 7) "This lead the model's accuracy to rise to 85%." (unchanged text; the
#    spell-check markers around "to" and the surrounding runs collapse together)
# ---------------------------------------------------------------------------
$r9 = $d.Content
$found9 = $r9.Find.Execute("rise to 85%.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "step7 found: $found9"
$r9.Text = "rise to 85%.ZZZTMP"
$r9b = $d.Content
$found9b = $r9b.Find.Execute("rise to 85%.ZZZTMP", $true, $false, $false, $false, $false, $true, 1, $false, "rise to 85%.", 2)
Write-Host "step7b found: $found9b"

# ---------------------------------------------------------------------------
# 8) "classifier still had good estimator" -> "classifier still had a good estimator"
# ---------------------------------------------------------------------------
$r10 = $d.Content
$found10 = $r10.Find.Execute("classifier still had good estimator", $true, $false, $false, $false, $false, $true, 1, $false, "classifier still had a good estimator", 2)
Write-Host "step8 found: $found10"
